$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns always carry their original text representation
# (avoids Excel auto-converting numeric-looking strings and dropping trailing zeros)
$ws.Range('D2').Value = '20.454.20'
$ws.Range('E2').Value = '  -7.31%  '
$ws.Range('D3').Value = '1.444.52'
$ws.Range('E3').Value = '  -7.29%  '
$ws.Range('D4').Value = '1.015'
$ws.Range('E4').Value = '  +1.54%  '
$ws.Range('D5').Value = '1.013'
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('D6').Value = '276.55'
$ws.Range('E6').Value = '  -5.17%  '
$ws.Range('D7').Value = '0.3705'
$ws.Range('E7').Value = '  -6.40%  '
$ws.Range('D8').Value = '0.3075'
$ws.Range('E8').Value = '  -5.15%  '
$ws.Range('D9').Value = '41.15'
$ws.Range('E9').Value = '  -7.55%  '
$ws.Range('E10').Value = '  -7.24%  '
$ws.Range('D11').Value = '0.06617'
$ws.Range('E11').Value = '  -9.17%  '
$ws.Range('D12').Value = '1.015'
$ws.Range('E12').Value = '  +1.50%  '
$ws.Range('D13').Value = '5.406'
$ws.Range('E13').Value = '  -5.67%  '
$ws.Range('D14').Value = '17.24'
$ws.Range('E14').Value = '  -8.56%  '
$ws.Range('D15').Value = '1.455.98'
$ws.Range('E15').Value = '  -6.47%  '
$ws.Range('D16').Value = '6.156'
$ws.Range('E16').Value = '  -7.61%  '
$ws.Range('D17').Value = '0.00001017'
$ws.Range('E17').Value = '  -9.95%  '
$ws.Range('D18').Value = '0.06364'
$ws.Range('E18').Value = '  -3.43%  '
$ws.Range('D19').Value = '77.89'
$ws.Range('E19').Value = '  -7.18%  '
$ws.Range('D20').Value = '1.015'
$ws.Range('E20').Value = '  +1.62%  '
$ws.Range('D21').Value = '5.736'
$ws.Range('E21').Value = '  -8.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.60'
$ws.Range('E22').Value = '  -6.62%  '
$ws.Range('D23').Value = '10.84'
$ws.Range('E23').Value = '  -4.62%  '
$ws.Range('D24').Value = '2.333'
$ws.Range('E24').Value = '  -1.77%  '
$ws.Range('D25').Value = '20.460.62'
$ws.Range('E25').Value = '  -7.34%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = '2.246'
$ws.Range('E26').Value = '  -7.71%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '142.70'
$ws.Range('E27').Value = '  -3.87%  '
$ws.Range('D28').Value = '17.18'
$ws.Range('E28').Value = '  -8.01%  '
$ws.Range('B29').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C29').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D29').Value = '1.625.22'
$ws.Range('E29').Value = '  -6.16%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '109.67'
$ws.Range('E30').Value = '  -8.38%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '0.9157'
$ws.Range('E31').Value = '  -7.81%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '5.495'
$ws.Range('E32').Value = '  -7.10%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').Value = '3.441'
$ws.Range('E33').Value = '  -29.31%  '
$ws.Range('D34').Value = '0.07724'
$ws.Range('E34').Value = '  -7.40%  '
$ws.Range('D35').Value = '8.291'
$ws.Range('E35').Value = '  -9.74%  '
$ws.Range('D36').Value = '1.435'
$ws.Range('E36').Value = '  -10.53%  '
$ws.Range('D37').Value = '1.013'
$ws.Range('E37').Value = '  +1.43%  '
$ws.Range('E38').Value = '  +1.07%  '
$ws.Range('D39').Value = '4.782'
$ws.Range('E39').Value = '  -7.05%  '
$ws.Range('D40').Value = '0.05581'
$ws.Range('E40').Value = '  -7.42%  '
$ws.Range('D41').Value = '0.02046'
$ws.Range('E41').Value = '  -10.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1910'
$ws.Range('E42').Value = '  -6.97%  '
$ws.Range('D43').Value = '1.115'
$ws.Range('E43').Value = '  -7.69%  '
$ws.Range('B44').Value = 'PancakeSwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D44').Value = '3.586'
$ws.Range('E44').Value = '  -4.55%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = '0.5337'
$ws.Range('E45').Value = '  -8.65%  '
$ws.Range('D46').Value = '12.16'
$ws.Range('E46').Value = '  -7.67%  '
$ws.Range('D47').Value = '0.5155'
$ws.Range('E47').Value = '  -8.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.770'
$ws.Range('E48').Value = '  -7.17%  '
$ws.Range('D49').Value = '109.28'
$ws.Range('E49').Value = '  -8.03%  '
$ws.Range('D50').Value = '1.066'
$ws.Range('E50').Value = '  -6.79%  '
$ws.Range('D51').Value = '0.06326'
$ws.Range('E51').Value = '  -7.33%  '
